# Weekly data refresh: a new price-report row for the week of 2023-03-31
# (serial date 45016) is inserted at row 191, pushing the existing
# historical rows 191:277 down to 192:278 (dimension grows from R277 to
# R278). The new row carries the same market/category metadata as its
# neighbours, with its own Volumen/Precio/Unidad values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 191, shifting rows 191-277 down to 192-278.
$ws.Rows("191:191").Insert()

# Populate the newly inserted row with the new week's data.
$ws.Range('A191').Value = 10
$ws.Range('B191').Value = 'Vega Modelo de Temuco'
$ws.Range('C191').Value = 'La Araucanía'
$ws.Range('D191').Value = 45016
$ws.Range('E191').Value = 9
$ws.Range('F191').Value = 100112005
$ws.Range('G191').Value = 'Puerro'
$ws.Range('H191').Value = 'Azul de Maquehue'
$ws.Range('I191').Value = 'Primera'
$ws.Range('J191').Value = 30
$ws.Range('K191').Value = 12000
$ws.Range('L191').Value = 12000
$ws.Range('M191').Value = 12000
$ws.Range('N191').Value = '$/docena de paquetes'
$ws.Range('O191').Value = 'Provincia de Cautín'
$ws.Range('P191').Value = 1000
$ws.Range('Q191').Value = 12
$ws.Range('R191').Value = 'Hortaliza'
